$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells AC1:AE1 ("Wins", "Losses", "Ties"), reusing the
# existing header formatting (bold / centered / bordered) by copying an
# existing header cell's format+value into place first, then overwriting
# the value - this avoids creating brand-new style/font records.
$ws.Range("AA1").Copy($ws.Range("AC1"))
$ws.Range("AA1").Copy($ws.Range("AD1"))
$ws.Range("AA1").Copy($ws.Range("AE1"))

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row.
$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 82  # AC -> Wins
    $ws.Cells.Item($r, 30).Value = 80  # AD -> Losses
    $ws.Cells.Item($r, 31).Value = 0   # AE -> Ties
}
